$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell's formatting (bold, border, centered) onto the
# two new header cells before filling in their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header cells I1, J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells I2:J4
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 5
